$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.362.82"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.936.98"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'0.7689"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.12%  "
$ws.Range("D6").Value = "'245.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.3202"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.46%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'27.81"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "'0.07027"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("D11").Value = "'0.7821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("D12").Value = "'0.08026"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "1.934.41"
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("D14").Value = "'5.353"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").Value = "'94.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "'14.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.77%  "
$ws.Range("D17").Value = "30.368.05"
$ws.Range("D18").Value = "'255.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'0.000007949"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.22%  "
$ws.Range("D20").Value = "'5.770"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "2.195.35"
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "'6.726"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.59%  "
$ws.Range("D25").Value = "'9.531"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "'164.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").Value = "'0.1345"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").Value = "'2.271"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.52%  "
$ws.Range("D30").Value = "'1.370"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").Value = "'4.408"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "'4.125"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "'0.7475"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'2.783"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'0.01955"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "
$ws.Range("D39").Value = "'2.812"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").Value = "'78.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D41").Value = "'6.407"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").Value = "'0.4497"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("D44").Value = "'1.003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "'0.8346"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'101.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "'9.788"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'7.497"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "'979.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.22%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("E51").Value = "  -1.60%  "
